$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.984.51'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.86%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.319.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -4.09%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.79%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.572'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.59%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.317.98'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.102'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.56'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.20%  '
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.336'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.738.03'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.935.96'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.80%  '
$ws.Range("E17").Value = '  -2.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.368.79'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '316.62'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.57'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.172'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.18%  '
$ws.Range("E28").Value = '  -6.96%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.81'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0731'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.12'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.82'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.384'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.18%  '
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.82'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.27'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '38.19'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.52'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '304.17'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '141.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.45'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0951'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0504'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.562'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.69'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0216'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.18%  '
$ws.Range("E51").Value = '  -0.18%  '
